$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(5)

# "TextBox 7" (Persona_occupation) was widened / nudged very slightly left:
#   off  x: 2061741 -> 2061740  (EMU)
#   ext cx: 1492843 -> 2178519  (EMU)
# COM exposes Left/Width in points (1 pt = 12700 EMU); the EMU<->pt
# round-trip through the host's float storage can land a unit off, so the
# literals below are nudged (within the same EMU bucket) to round-trip to
# the exact target EMU values.
$sh.Left  = 162.34173228346458
$sh.Width = 171.536974
